{"js": "// Replace the 100 arithmetic-expression answers that live in the single\n// 20-row x 5-column table of the document. Each table cell holds exactly\n// one paragraph with one run; we overwrite the run's visible text only\n// (in place, via a Range replace) so that every other property of the\n// run/paragraph (fonts, size, justification, etc.) is left untouched,\n// matching the diff which only touches the <w:t> contents.\n\n// New value for every cell, in row-major order (row 0 col 0..4, row 1 col\n// 0..4, ...), taken 1:1 from the authoritative diff.\nconst newValues = [\n  \"74-69=5\", \"2+0=2\", \"2+86=88\", \"14+36=50\", \"88-49=39\",\n  \"1+23=24\", \"17+79=96\", \"40-8=32\", \"65-20=45\", \"22+3=25\",\n  \"95-47=48\", \"20+62=82\", \"98-28=70\", \"61+13=74\", \"7+59=66\",\n  \"26-4=22\", \"25+4=29\", \"93-62=31\", \"7+87=94\", \"80-15=65\",\n  \"16+72=88\", \"13+48=61\", \"98-77=21\", \"57-8=49\", \"54-52=2\",\n  \"43+31=74\", \"60+1=61\", \"1+36=37\", \"9+38=47\", \"57-5=52\",\n  \"24+5=29\", \"54+21=75\", \"30-12=18\", \"97-41=56\", \"88-76=12\",\n  \"52+38=90\", \"4+14=18\", \"81+9=90\", \"62-40=22\", \"70-19=51\",\n  \"14-12=2\", \"36+1=37\", \"10+7=17\", \"97-16=81\", \"71+27=98\",\n  \"99-49=50\", \"21+53=74\", \"12+38=50\", \"93-42=51\", \"34-24=10\",\n  \"56+28=84\", \"97-73=24\", \"39-9=30\", \"79-43=36\", \"63-9=54\",\n  \"76-59=17\", \"17+8=25\", \"33-16=17\", \"7+46=53\", \"48+49=97\",\n  \"77-1=76\", \"24+49=73\", \"73-9=64\", \"60-10=50\", \"58-47=11\",\n  \"33+4=37\", \"82-39=43\", \"27+71=98\", \"69+30=99\", \"8+29=37\",\n  \"77-28=49\", \"96+1=97\", \"38-37=1\", \"55-32=23\", \"96-36=60\",\n  \"5+62=67\", \"1+17=18\", \"17+73=90\", \"15+76=91\", \"94-90=4\",\n  \"19-8=11\", \"31+30=61\", \"9+67=76\", \"13-12=1\", \"4+32=36\",\n  \"85-27=58\", \"26+62=88\", \"43-5=38\", \"61+22=83\", \"97-76=21\",\n  \"89+10=99\", \"45+51=96\", \"55+14=69\", \"47+44=91\", \"86-27=59\",\n  \"97-75=22\", \"65+19=84\", \"2+41=43\", \"37+52=89\", \"18+1=19\",\n];\n\nconst ROWS = 20;\nconst COLS = 5;\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (let r = 0; r < ROWS; r++) {\n  for (let c = 0; c < COLS; c++) {\n    const idx = r * COLS + c;\n    const cell = table.getCell(r, c);\n    const para = cell.body.paragraphs.getFirst();\n    const range = para.getRange();\n    range.insertText(newValues[idx], Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 100 arithmetic-expression answers that live in the single\n# 20-row x 5-column table of the document. Each table cell holds exactly\n# one paragraph with one run; setting Cell.Range.Text overwrites the\n# run's visible text only (Word keeps the run's rPr / paragraph's pPr in\n# place), matching the diff which only touches the <w:t> contents.\n\n$d = $word.ActiveDocument\n\n# New value for every cell, in row-major order (row 1 col 1..5, row 2 col\n# 1..5, ...), taken 1:1 from the authoritative diff.\n$newValues = @(\n  \"74-69=5\", \"2+0=2\", \"2+86=88\", \"14+36=50\", \"88-49=39\",\n  \"1+23=24\", \"17+79=96\", \"40-8=32\", \"65-20=45\", \"22+3=25\",\n  \"95-47=48\", \"20+62=82\", \"98-28=70\", \"61+13=74\", \"7+59=66\",\n  \"26-4=22\", \"25+4=29\", \"93-62=31\", \"7+87=94\", \"80-15=65\",\n  \"16+72=88\", \"13+48=61\", \"98-77=21\", \"57-8=49\", \"54-52=2\",\n  \"43+31=74\", \"60+1=61\", \"1+36=37\", \"9+38=47\", \"57-5=52\",\n  \"24+5=29\", \"54+21=75\", \"30-12=18\", \"97-41=56\", \"88-76=12\",\n  \"52+38=90\", \"4+14=18\", \"81+9=90\", \"62-40=22\", \"70-19=51\",\n  \"14-12=2\", \"36+1=37\", \"10+7=17\", \"97-16=81\", \"71+27=98\",\n  \"99-49=50\", \"21+53=74\", \"12+38=50\", \"93-42=51\", \"34-24=10\",\n  \"56+28=84\", \"97-73=24\", \"39-9=30\", \"79-43=36\", \"63-9=54\",\n  \"76-59=17\", \"17+8=25\", \"33-16=17\", \"7+46=53\", \"48+49=97\",\n  \"77-1=76\", \"24+49=73\", \"73-9=64\", \"60-10=50\", \"58-47=11\",\n  \"33+4=37\", \"82-39=43\", \"27+71=98\", \"69+30=99\", \"8+29=37\",\n  \"77-28=49\", \"96+1=97\", \"38-37=1\", \"55-32=23\", \"96-36=60\",\n  \"5+62=67\", \"1+17=18\", \"17+73=90\", \"15+76=91\", \"94-90=4\",\n  \"19-8=11\", \"31+30=61\", \"9+67=76\", \"13-12=1\", \"4+32=36\",\n  \"85-27=58\", \"26+62=88\", \"43-5=38\", \"61+22=83\", \"97-76=21\",\n  \"89+10=99\", \"45+51=96\", \"55+14=69\", \"47+44=91\", \"86-27=59\",\n  \"97-75=22\", \"65+19=84\", \"2+41=43\", \"37+52=89\", \"18+1=19\"\n)\n\n$table = $d.Tables.Item(1)\n$ROWS = 20\n$COLS = 5\n\nfor ($r = 1; $r -le $ROWS; $r++) {\n    for ($c = 1; $c -le $COLS; $c++) {\n        $idx = ($r - 1) * $COLS + ($c - 1)\n        $cell = $table.Cell($r, $c)\n        $cell.Range.Text = $newValues[$idx]\n    }\n}\n"}
